$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 35, shifting existing rows 35..157 down to 36..158
$ws.Rows("35:35").Insert([Microsoft.Office.Interop.Excel.XlInsertShiftDirection]::xlShiftDown)

# Populate the newly inserted row 35 with the new record's data
$ws.Range("A35").Value = 5
$ws.Range("B35").Value = 'Macroferia Regional de Talca'
$ws.Range("C35").Value = 'Maule'
$ws.Range("D35").Value = 44672
$ws.Range("E35").Value = 7
$ws.Range("F35").Value = 100112017
$ws.Range("G35").Value = 'Apio'
$ws.Range("H35").Value = 'Americana (o)'
$ws.Range("I35").Value = 'Primera'
$ws.Range("J35").Value = 500
$ws.Range("K35").Value = 9000
$ws.Range("L35").Value = 9000
$ws.Range("M35").Value = 9000
$ws.Range("N35").Value = '$/docena de matas'
$ws.Range("O35").Value = 'Provincia del Elquí'
$ws.Range("P35").Value = 1500
$ws.Range("Q35").Value = 6
$ws.Range("R35").Value = 'Hortaliza'

# Match the date-number format used by the rest of column D
$ws.Range("D35").NumberFormat = $ws.Range("D36").NumberFormat
